$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item("总计")

# Duplicate the existing "总计" sheet so the copy becomes the NEW "总计"
# (keeps identical sheetPr / pageSetup / header style conventions)
$total.Copy([System.Reflection.Missing]::Value, $total)
$newTotal = $wb.Worksheets.Item($wb.Worksheets.Count)

# ---- Shift the old rows down on the copy and insert the new 2022-Q1 summary row ----
$newTotal.Range("A4").Value = 2
$newTotal.Range("B4").Value = "2021-Q3"
$newTotal.Range("C4").Value = 17
$newTotal.Range("D4").Value = 5.19
$newTotal.Range("A3").Copy()
$newTotal.Range("A4").PasteSpecial(-4122)

$newTotal.Range("A3").Value = 1
$newTotal.Range("B3").Value = "2021-Q4"
$newTotal.Range("C3").Value = 24
$newTotal.Range("D3").Value = 4.86

$newTotal.Range("A2").Value = 0
$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 14
$newTotal.Range("D2").Value = 0.38

# ---- Repurpose the ORIGINAL "总计" worksheet object into the new "2022-Q1" fund-holdings sheet ----
# (rename the original out of the way FIRST so the copy can take the "总计" name cleanly)
$q1 = $total
$q1.Name = "2022-Q1"
$newTotal.Name = "总计"

# Keep these numeric-looking strings as literal text (matches source data formatting).
# (Column C holds fund names, which are never numeric-looking, so it is left on General;
#  G15 is a real 0 value, not text, so it is excluded from column G's text range below.)
$q1.Range("B2:B15").NumberFormat = "@"
$q1.Range("D2:D15").NumberFormat = "@"
$q1.Range("E2:E15").NumberFormat = "@"
$q1.Range("F2:F15").NumberFormat = "@"
$q1.Range("G2:G14").NumberFormat = "@"

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Headers need the same bold/border/center "s=2" style as B1:D1 (inherited from the old sheet)
$q1.Range("D1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "009381"
$q1.Range("C2").Value = "汇安核心资产混合A"
$q1.Range("D2").Value = "4.50"
$q1.Range("E2").Value = "94.22"
$q1.Range("F2").Value = "3.11"
$q1.Range("G2").Value = "0.1400"
$q1.Range("H2").Value = 10

$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "010558"
$q1.Range("C3").Value = "汇安鑫利优选混合A"
$q1.Range("D3").Value = "2.21"
$q1.Range("E3").Value = "92.83"
$q1.Range("F3").Value = "3.02"
$q1.Range("G3").Value = "0.0667"
$q1.Range("H3").Value = 9

$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "003889"
$q1.Range("C4").Value = "汇安丰泽灵活配置混合A"
$q1.Range("D4").Value = "1.36"
$q1.Range("E4").Value = "93.28"
$q1.Range("F4").Value = "3.00"
$q1.Range("G4").Value = "0.0408"
$q1.Range("H4").Value = 9

$q1.Range("A5").Value = 3
$q1.Range("B5").Value = "004558"
$q1.Range("C5").Value = "汇安丰裕灵活配置混合A"
$q1.Range("D5").Value = "0.99"
$q1.Range("E5").Value = "83.41"
$q1.Range("F5").Value = "3.05"
$q1.Range("G5").Value = "0.0302"
$q1.Range("H5").Value = 7

$q1.Range("A6").Value = 4
$q1.Range("B6").Value = "004560"
$q1.Range("C6").Value = "汇安丰益灵活配置混合A"
$q1.Range("D6").Value = "2.04"
$q1.Range("E6").Value = "30.34"
$q1.Range("F6").Value = "1.15"
$q1.Range("G6").Value = "0.0235"
$q1.Range("H6").Value = 9

$q1.Range("A7").Value = 5
$q1.Range("B7").Value = "010559"
$q1.Range("C7").Value = "汇安鑫利优选混合C"
$q1.Range("D7").Value = "0.73"
$q1.Range("E7").Value = "92.83"
$q1.Range("F7").Value = "3.02"
$q1.Range("G7").Value = "0.0220"
$q1.Range("H7").Value = 9

$q1.Range("A8").Value = 6
$q1.Range("B8").Value = "005599"
$q1.Range("C8").Value = "汇安量化优选灵活配置混合A"
$q1.Range("D8").Value = "0.85"
$q1.Range("E8").Value = "93.69"
$q1.Range("F8").Value = "2.22"
$q1.Range("G8").Value = "0.0189"
$q1.Range("H8").Value = 9

$q1.Range("A9").Value = 7
$q1.Range("B9").Value = "003890"
$q1.Range("C9").Value = "汇安丰泽灵活配置混合C"
$q1.Range("D9").Value = "0.51"
$q1.Range("E9").Value = "93.28"
$q1.Range("F9").Value = "3.00"
$q1.Range("G9").Value = "0.0153"
$q1.Range("H9").Value = 9

$q1.Range("A10").Value = 8
$q1.Range("B10").Value = "007775"
$q1.Range("C10").Value = "汇安量化先锋混合A"
$q1.Range("D10").Value = "0.38"
$q1.Range("E10").Value = "94.51"
$q1.Range("F10").Value = "3.20"
$q1.Range("G10").Value = "0.0122"
$q1.Range("H10").Value = 6

$q1.Range("A11").Value = 9
$q1.Range("B11").Value = "007776"
$q1.Range("C11").Value = "汇安量化先锋混合C"
$q1.Range("D11").Value = "0.11"
$q1.Range("E11").Value = "94.51"
$q1.Range("F11").Value = "3.20"
$q1.Range("G11").Value = "0.0035"
$q1.Range("H11").Value = 6

$q1.Range("A12").Value = 10
$q1.Range("B12").Value = "009382"
$q1.Range("C12").Value = "汇安核心资产混合C"
$q1.Range("D12").Value = "0.08"
$q1.Range("E12").Value = "94.22"
$q1.Range("F12").Value = "3.11"
$q1.Range("G12").Value = "0.0025"
$q1.Range("H12").Value = 10

$q1.Range("A13").Value = 11
$q1.Range("B13").Value = "005600"
$q1.Range("C13").Value = "汇安量化优选灵活配置混合C"
$q1.Range("D13").Value = "0.07"
$q1.Range("E13").Value = "93.69"
$q1.Range("F13").Value = "2.22"
$q1.Range("G13").Value = "0.0016"
$q1.Range("H13").Value = 9

$q1.Range("A14").Value = 12
$q1.Range("B14").Value = "004559"
$q1.Range("C14").Value = "汇安丰裕灵活配置混合C"
$q1.Range("D14").Value = "0.01"
$q1.Range("E14").Value = "83.41"
$q1.Range("F14").Value = "3.05"
$q1.Range("G14").Value = "0.0003"
$q1.Range("H14").Value = 7

$q1.Range("A15").Value = 13
$q1.Range("B15").Value = "004561"
$q1.Range("C15").Value = "汇安丰益灵活配置混合C"
$q1.Range("D15").Value = "0.00"
$q1.Range("E15").Value = "30.34"
$q1.Range("F15").Value = "1.15"
$q1.Range("G15").Value = 0
$q1.Range("H15").Value = 9

# Re-apply the bold/border "s=2" look to column A (rows 2-15), matching the sibling sheets
$wb.Worksheets.Item("2021-Q4").Range("A2").Copy()
$q1.Range("A2:A15").PasteSpecial(-4122)

# Move the new "总计" sheet to the very end (after "2022-Q1")
$newTotal.Move([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))

# Restore the originally active sheet/tab
$wb.Worksheets.Item("2021-Q3").Activate()
Write-Host "done"